# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 6 de Mayo de 2020 a las 14:03"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1238083
$ws.Range("C4").Value = 450
$ws.Range("D4").Value = 201011
$ws.Range("E4").Value = 964787
$ws.Range("G4").Value = 14
$ws.Range("H4").Value = 72285

# Row 5 - España
$ws.Range("B5").Value = 253682
$ws.Range("C5").Value = 3121
$ws.Range("D5").Value = 159359
$ws.Range("E5").Value = 68466
$ws.Range("F5").Value = 2075
$ws.Range("G5").Value = 244
$ws.Range("H5").Value = 25857

# Row 12 - Brasil
$ws.Range("B12").Value = 116299
$ws.Range("C12").Value = 1584
$ws.Range("E12").Value = 60112
$ws.Range("G12").Value = 45
$ws.Range("H12").Value = 7966

# Row 19 - Paises Bajos
$ws.Range("B19").Value = 41319
$ws.Range("C19").Value = 232
$ws.Range("E19").Value = 35865
$ws.Range("G19").Value = 36
$ws.Range("H19").Value = 5204

# Row 55 - Finlandia
$ws.Range("E55").Value = 1821
$ws.Range("F55").Value = 44
$ws.Range("G55").Value = 6
$ws.Range("H55").Value = 252

# Row 60 - Kazajistan
$ws.Range("D60").Value = 1387
$ws.Range("E60").Value = 2881
$ws.Range("F60").Value = 31
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 30

# Row 76 - Bosnia y Herzegovina
$ws.Range("B76").Value = 1987
$ws.Range("C76").Value = 41
$ws.Range("D76").Value = 928
$ws.Range("E76").Value = 973
$ws.Range("G76").Value = 7
$ws.Range("H76").Value = 86

# Row 83 - Republica de Macedonia
$ws.Range("B83").Value = 1539
$ws.Range("C83").Value = 13
$ws.Range("D83").Value = 1057
$ws.Range("E83").Value = 394
$ws.Range("G83").Value = 2
$ws.Range("H83").Value = 88

# Rows 87-89: Senegal's updated figures push it above Eslovaquia and Lituania
# Row 87 becomes Senegal (was Eslovaquia)
$ws.Range("A87").Value = "Senegal"
$ws.Range("B87").Value = 1433
$ws.Range("C87").Value = 104
$ws.Range("D87").Value = 493
$ws.Range("E87").Value = 928
$ws.Range("F87").Value = 6
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 12

# Row 88 becomes Eslovaquia (was Lituania), keeping former row-87 figures
$ws.Range("A88").Value = "Eslovaquia"
$ws.Range("B88").Value = 1429
$ws.Range("C88").Value = 8
$ws.Range("D88").Value = 762
$ws.Range("E88").Value = 642
$ws.Range("F88").Value = 4
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 25

# Row 89 becomes Lituania (was Senegal), keeping former row-88 figures
$ws.Range("A89").Value = "Lituania"
$ws.Range("B89").Value = 1428
$ws.Range("C89").Value = 5
$ws.Range("D89").Value = 718
$ws.Range("E89").Value = 662
$ws.Range("F89").Value = 17
$ws.Range("G89").Value = 2
$ws.Range("H89").Value = 48

# Row 117 - Malta
$ws.Range("B117").Value = 484
$ws.Range("C117").Value = 2
$ws.Range("D117").Value = 407
$ws.Range("E117").Value = 72
$ws.Range("F117").Value = 0

# Rows 205-206: Seychelles overtakes Montserrat
# Row 205 becomes Seychelles (was Montserrat)
$ws.Range("A205").Value = "Seychelles"
$ws.Range("D205").Value = 8
$ws.Range("F205").Value = 0
$ws.Range("H205").Value = 0

# Row 206 becomes Montserrat (was Seychelles)
$ws.Range("A206").Value = "Montserrat"
$ws.Range("D206").Value = 7
$ws.Range("F206").Value = 1
$ws.Range("H206").Value = 1
